$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Y33: was stored as text " 7.6355", should be the numeric value 7.6355
$ws.Range("Y33").Value = 7.6355000000000004

# Column AA ("note"): fill every row that doesn't already have a note
# with a copy of the battle name from column A (used so every circle on
# the map has an on-click popup label even when there's no special note).
for ($r = 2; $r -le 39; $r++) {
    $noteCell = $ws.Cells.Item($r, 27)
    if ($noteCell.Value2 -eq $null) {
        $ws.Cells.Item($r, 27).Value = $ws.Cells.Item($r, 1).Value2
    }
}

# Nudge AA39's font so it carries an explicit (applied) font style, matching
# the workbook author's formatting pass on the last row.
$ws.Range("AA39").Font.ThemeColor = 1

# Restore the view state (selected cell) from the session that produced the edit.
$ws.Range("Y34").Select() | Out-Null
